$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.172.56"
$ws.Range("E2").Value = "  -2.25%  "

$ws.Range("D3").Value = "1.852.59"
$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.41"
$ws.Range("E5").Value = "  -1.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6823"
$ws.Range("E6").Value = "  -6.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07708"
$ws.Range("E8").Value = "  +1.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3034"
$ws.Range("E9").Value = "  -3.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.11"
$ws.Range("E10").Value = "  -6.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08162"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").Value = "1.899.28"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7224"
$ws.Range("E13").Value = "  -3.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.191"
$ws.Range("E14").Value = "  -3.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.23"
$ws.Range("E15").Value = "  -3.78%  "

$ws.Range("D16").Value = "29.163.10"
$ws.Range("E16").Value = "  -2.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007796"
$ws.Range("E17").Value = "  -1.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.720"
$ws.Range("E18").Value = "  -4.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.13"
$ws.Range("E19").Value = "  -2.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "233.71"
$ws.Range("E20").Value = "  -5.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "2.099.95"
$ws.Range("E22").Value = "  -1.14%  "

$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.413"
$ws.Range("E24").Value = "  -4.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.67"
$ws.Range("E25").Value = "  -1.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.945"
$ws.Range("E26").Value = "  -3.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1427"
$ws.Range("E27").Value = "  -6.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.03"
$ws.Range("E28").Value = "  -3.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.959"
$ws.Range("E29").Value = "  -2.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.399"
$ws.Range("E30").Value = "  -3.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.514"
$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("E32").Value = "  -2.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.005"
$ws.Range("E33").Value = "  -4.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05180"
$ws.Range("E34").Value = "  -4.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.178"
$ws.Range("E35").Value = "  -4.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7027"
$ws.Range("E36").Value = "  -5.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.025"
$ws.Range("E37").Value = "  +1.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.662"
$ws.Range("E38").Value = "  -1.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01845"
$ws.Range("E39").Value = "  -4.48%  "

$ws.Range("E40").Value = "  -2.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9131"
$ws.Range("E41").Value = "  +2.76%  "

$ws.Range("D42").Value = "1.104.13"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.972"
$ws.Range("E43").Value = "  -0.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4267"
$ws.Range("E44").Value = "  -4.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.22"
$ws.Range("E45").Value = "  -2.21%  "

$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.43"
$ws.Range("E47").Value = "  -1.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.758"
$ws.Range("E48").Value = "  -3.67%  "

$ws.Range("D49").Value = "1.997.21"
$ws.Range("E49").Value = "  -1.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.139"
$ws.Range("E50").Value = "  -5.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.883"
$ws.Range("E51").Value = "  -8.15%  "
